$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (G2=5489) in ALC
$ws.Range("H2").Value2 = 406.53845
$ws.Range("I2").Value2 = 380.1
$ws.Range("J2").Value2 = 494.66666
$ws.Range("K2").Value2 = 380.1
$ws.Range("L2").Value2 = 494.66666
$ws.Range("M2").Value2 = -267.1
$ws.Range("N2").Value2 = -720.66666

# Row 6 (G6=4564) in ALC
$ws.Range("H6").Value2 = 704.55554
$ws.Range("I6").Value2 = 740.1667
$ws.Range("J6").Value2 = 633.3333
$ws.Range("K6").Value2 = 2220.5001
$ws.Range("L6").Value2 = 1899.9999
$ws.Range("M6").Value2 = -2108.5001
$ws.Range("N6").Value2 = -2123.9999

# Row 12 (G12=5515) in ALC
$ws.Range("H12").Value2 = 168.14285
$ws.Range("I12").Value2 = 162.5
$ws.Range("J12").Value2 = 202
$ws.Range("K12").Value2 = 162.5
$ws.Range("L12").Value2 = 202
$ws.Range("M12").Value2 = 7.5
$ws.Range("N12").Value2 = -542

# Row 21 (G21=2149) in ALC
$ws.Range("H21").Value2 = 18666.666
$ws.Range("I21").Value2 = 0
$ws.Range("J21").Value2 = 18666.666
$ws.Range("K21").Value2 = 0
$ws.Range("L21").Value2 = 18666.666
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value2 = -19602.666

# Row 23 (G23=2149) in ALC
$ws.Range("H23").Value2 = 18666.666
$ws.Range("I23").Value2 = 0
$ws.Range("J23").Value2 = 18666.666
$ws.Range("K23").Value2 = 0
$ws.Range("L23").Value2 = 18666.666
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value2 = -19134.666

# Row 29 (G29=4575) in ALC
$ws.Range("H29").Value2 = 225
$ws.Range("I29").Value2 = 225
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 675
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = -394
$ws.Range("N29").ClearContents()

# Row 52 (G52=4567) in ALC
$ws.Range("H52").Value2 = 1533.3334
$ws.Range("I52").Value2 = 1000
$ws.Range("J52").Value2 = 1800
$ws.Range("K52").Value2 = 3000
$ws.Range("L52").Value2 = 5400
$ws.Range("M52").Value2 = -2840
$ws.Range("N52").Value2 = -5720

# Row 62 (G62=27781) in ALC
$ws.Range("H62").Value2 = 5000
$ws.Range("I62").Value2 = 8300
$ws.Range("J62").Value2 = 1700
$ws.Range("K62").Value2 = 8300
$ws.Range("L62").Value2 = 1700
$ws.Range("M62").Value2 = -7676
$ws.Range("N62").Value2 = -2948

# Row 65 (G65=27781) in ALC
$ws.Range("H65").Value2 = 5000
$ws.Range("I65").Value2 = 8300
$ws.Range("J65").Value2 = 1700
$ws.Range("K65").Value2 = 41500
$ws.Range("L65").Value2 = 8500
$ws.Range("M65").Value2 = -38380
$ws.Range("N65").Value2 = -14740

# Row 94 (G94=19905) in ALC
$ws.Range("H94").Value2 = 4222.222
$ws.Range("I94").Value2 = 3333.3333
$ws.Range("K94").Value2 = 3333.3333
$ws.Range("M94").Value2 = -2882.3333

# Row 137 (G137=44013) in ALC
$ws.Range("H137").Value2 = 15387269
$ws.Range("I137").Value2 = 783.1667
$ws.Range("J137").Value2 = 28575686
$ws.Range("K137").Value2 = 2349.5001
$ws.Range("L137").Value2 = 85727058
$ws.Range("M137").Value2 = 200.4998999999998
$ws.Range("N137").Value2 = -85732158

# Row 141 (G141=44161) in ALC
$ws.Range("H141").Value2 = 3320.8823
$ws.Range("I141").Value2 = 1175
$ws.Range("J141").Value2 = 4491.364
$ws.Range("K141").Value2 = 3525
$ws.Range("L141").Value2 = 13474.092
$ws.Range("M141").Value2 = 1655
$ws.Range("N141").Value2 = -23834.092

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G2=27713) in ARM
$ws.Range("H2").Value2 = 1680.9231
$ws.Range("I2").Value2 = 1329.875
$ws.Range("J2").Value2 = 2242.6
$ws.Range("K2").Value2 = 1329.875
$ws.Range("L2").Value2 = 2242.6
$ws.Range("M2").Value2 = -1216.875
$ws.Range("N2").Value2 = -2468.6

# Row 45 (G45=27714) in ARM
$ws.Range("H45").Value2 = 2165897.2
$ws.Range("I45").Value2 = 2842116.8
$ws.Range("K45").Value2 = 2842116.8
$ws.Range("M45").Value2 = -2841739.8

# Row 70 (G70=19555) in ARM
$ws.Range("H70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("N70").ClearContents()

# Row 73 (G73=19555) in ARM
$ws.Range("H73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("N73").ClearContents()

# Row 102 (G102=19945) in ARM
$ws.Range("H102").Value2 = 2500
$ws.Range("I102").Value2 = 0
$ws.Range("J102").Value2 = 2500
$ws.Range("K102").Value2 = 0
$ws.Range("L102").Value2 = 2500
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value2 = -5744

# Row 116 (G116=27713) in ARM
$ws.Range("H116").Value2 = 1680.9231
$ws.Range("I116").Value2 = 1329.875
$ws.Range("J116").Value2 = 2242.6
$ws.Range("K116").Value2 = 1329.875
$ws.Range("L116").Value2 = 2242.6
$ws.Range("M116").Value2 = 964.125
$ws.Range("N116").Value2 = -6830.6

# Row 122 (G122=36168) in ARM
$ws.Range("H122").Value2 = 9228.625
$ws.Range("I122").Value2 = 9725.736999999999
$ws.Range("J122").Value2 = 7339.6
$ws.Range("K122").Value2 = 29177.211
$ws.Range("L122").Value2 = 22018.8
$ws.Range("M122").Value2 = -26727.211
$ws.Range("N122").Value2 = -26918.8

# Row 132 (G132=43997) in ARM
$ws.Range("H132").Value2 = 1844.5454
$ws.Range("I132").Value2 = 1661.9
$ws.Range("K132").Value2 = 4985.700000000001
$ws.Range("M132").Value2 = -2455.700000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G3=27713) in BSM
$ws.Range("H3").Value2 = 1680.9231
$ws.Range("I3").Value2 = 1329.875
$ws.Range("J3").Value2 = 2242.6
$ws.Range("K3").Value2 = 1329.875
$ws.Range("L3").Value2 = 2242.6
$ws.Range("M3").Value2 = -1215.875
$ws.Range("N3").Value2 = -2470.6

# Row 99 (G99=19943) in BSM
$ws.Range("H99").Value2 = 1477
$ws.Range("I99").Value2 = 1136.6666
$ws.Range("J99").Value2 = 2498
$ws.Range("K99").Value2 = 1136.6666
$ws.Range("L99").Value2 = 2498
$ws.Range("M99").Value2 = 361.3334
$ws.Range("N99").Value2 = -5494

# Row 134 (G134=43998) in BSM
$ws.Range("H134").Value2 = 3079.158
$ws.Range("I134").Value2 = 2852.32
$ws.Range("J134").Value2 = 3515.3845
$ws.Range("K134").Value2 = 8556.960000000001
$ws.Range("L134").Value2 = 10546.1535
$ws.Range("M134").Value2 = -6021.960000000001
$ws.Range("N134").Value2 = -15616.1535

$ws = $wb.Worksheets.Item("CRP")
# Row 122 (G122=36196) in CRP
$ws.Range("H122").Value2 = 2119.2
$ws.Range("I122").Value2 = 2260.6667
$ws.Range("J122").Value2 = 1907
$ws.Range("K122").Value2 = 6782.000100000001
$ws.Range("L122").Value2 = 5721
$ws.Range("M122").Value2 = -4332.000100000001
$ws.Range("N122").Value2 = -10621

# Row 132 (G132=44019) in CRP
$ws.Range("H132").Value2 = 19232578
$ws.Range("I132").Value2 = 21740594
$ws.Range("K132").Value2 = 65221782
$ws.Range("M132").Value2 = -65219252

$ws = $wb.Worksheets.Item("CUL")
# Row 36 (G36=4732) in CUL
$ws.Range("H36").Value2 = 900
$ws.Range("I36").Value2 = 900
$ws.Range("K36").Value2 = 2700
$ws.Range("M36").Value2 = -2531

# Row 63 (G63=12866) in CUL
$ws.Range("H63").Value2 = 9622.571
$ws.Range("I63").Value2 = 0
$ws.Range("J63").Value2 = 9622.571
$ws.Range("K63").Value2 = 0
$ws.Range("L63").Value2 = 28867.713
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value2 = -30365.713

# Row 64 (G64=12861) in CUL
$ws.Range("H64").Value2 = 5491.6665
$ws.Range("I64").Value2 = 750
$ws.Range("J64").Value2 = 6440
$ws.Range("K64").Value2 = 2250
$ws.Range("L64").Value2 = 19320
$ws.Range("M64").Value2 = -1980
$ws.Range("N64").Value2 = -19860

# Row 66 (G66=12866) in CUL
$ws.Range("H66").Value2 = 9622.571
$ws.Range("I66").Value2 = 0
$ws.Range("J66").Value2 = 9622.571
$ws.Range("K66").Value2 = 0
$ws.Range("L66").Value2 = 86603.139
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value2 = -94091.139

# Row 67 (G67=12861) in CUL
$ws.Range("H67").Value2 = 5491.6665
$ws.Range("I67").Value2 = 750
$ws.Range("J67").Value2 = 6440
$ws.Range("K67").Value2 = 2250
$ws.Range("L67").Value2 = 19320
$ws.Range("M67").Value2 = -1314
$ws.Range("N67").Value2 = -21192

# Row 114 (G114=27865) in CUL
$ws.Range("H114").Value2 = 1310.6
$ws.Range("I114").Value2 = 461.33334
$ws.Range("J114").Value2 = 1578.7894
$ws.Range("K114").Value2 = 1384.00002
$ws.Range("L114").Value2 = 4736.3682
$ws.Range("M114").Value2 = 1869.99998
$ws.Range("N114").Value2 = -11244.3682

# Row 132 (G132=43972) in CUL
$ws.Range("H132").Value2 = 1249.0714
$ws.Range("I132").Value2 = 765.2222
$ws.Range("J132").Value2 = 2120
$ws.Range("K132").Value2 = 6886.999800000001
$ws.Range("L132").Value2 = 19080
$ws.Range("M132").Value2 = -4356.999800000001
$ws.Range("N132").Value2 = -24140

$ws = $wb.Worksheets.Item("GSM")
# Row 118 (G118=26172) in GSM
$ws.Range("H118").Value2 = 14281.25
$ws.Range("J118").Value2 = 14281.25
$ws.Range("L118").Value2 = 14281.25
$ws.Range("N118").Value2 = -17595.25

# Row 122 (G122=36182) in GSM
$ws.Range("H122").Value2 = 3923663
$ws.Range("I122").Value2 = 5129736.5
$ws.Range("J122").Value2 = 3924.5
$ws.Range("K122").Value2 = 15389209.5
$ws.Range("L122").Value2 = 11773.5
$ws.Range("M122").Value2 = -15386759.5
$ws.Range("N122").Value2 = -16673.5

# Row 132 (G132=44008) in GSM
$ws.Range("H132").Value2 = 4551.3184
$ws.Range("I132").Value2 = 4706.3057
$ws.Range("J132").Value2 = 3853.875
$ws.Range("K132").Value2 = 14118.9171
$ws.Range("L132").Value2 = 11561.625
$ws.Range("M132").Value2 = -11588.9171
$ws.Range("N132").Value2 = -16621.625

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G22=5277) in LTW
$ws.Range("H22").Value2 = 1739.25
$ws.Range("I22").Value2 = 1300
$ws.Range("J22").Value2 = 1802
$ws.Range("K22").Value2 = 1300
$ws.Range("L22").Value2 = 1802
$ws.Range("M22").Value2 = -1005
$ws.Range("N22").Value2 = -2392

# Row 27 (G27=5277) in LTW
$ws.Range("H27").Value2 = 1739.25
$ws.Range("I27").Value2 = 1300
$ws.Range("J27").Value2 = 1802
$ws.Range("K27").Value2 = 1300
$ws.Range("L27").Value2 = 1802
$ws.Range("M27").Value2 = -1193
$ws.Range("N27").Value2 = -2016

# Row 46 (G46=5282) in LTW
$ws.Range("H46").Value2 = 886.55884
$ws.Range("I46").Value2 = 631.6
$ws.Range("J46").Value2 = 992.7917
$ws.Range("K46").Value2 = 631.6
$ws.Range("L46").Value2 = 992.7917
$ws.Range("M46").Value2 = -443.6
$ws.Range("N46").Value2 = -1368.7917

# Row 60 (G60=3904) in LTW
$ws.Range("H60").Value2 = 25745.5
$ws.Range("J60").Value2 = 25745.5
$ws.Range("L60").Value2 = 25745.5
$ws.Range("N60").Value2 = -26763.5

# Row 100 (G100=19995) in LTW
$ws.Range("H100").Value2 = 2598.25
$ws.Range("I100").Value2 = 2339.2
$ws.Range("J100").Value2 = 3030
$ws.Range("K100").Value2 = 2339.2
$ws.Range("L100").Value2 = 3030
$ws.Range("M100").Value2 = -1798.2
$ws.Range("N100").Value2 = -4112

# Row 122 (G122=36247) in LTW
$ws.Range("H122").Value2 = 5683.9736
$ws.Range("I122").Value2 = 5622.7085
$ws.Range("J122").Value2 = 5789
$ws.Range("K122").Value2 = 16868.1255
$ws.Range("L122").Value2 = 17367
$ws.Range("M122").Value2 = -14418.1255
$ws.Range("N122").Value2 = -22267

$ws = $wb.Worksheets.Item("WVR")
# Row 59 (G59=3201) in WVR
$ws.Range("H59").Value2 = 0
$ws.Range("J59").Value2 = 0
$ws.Range("L59").Value2 = 0
$ws.Range("N59").ClearContents()

# Row 132 (G132=44029) in WVR
$ws.Range("H132").Value2 = 1057.9615
$ws.Range("I132").Value2 = 895.7805
$ws.Range("J132").Value2 = 1662.4546
$ws.Range("K132").Value2 = 2687.3415
$ws.Range("L132").Value2 = 4987.3638
$ws.Range("M132").Value2 = -157.3415
$ws.Range("N132").Value2 = -10047.3638
